# Applies the odds updates for Jogos_da_Semana_FlashScore_2025-05-27.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.95

# Row 5
$ws.Range("G5").Value = 8.75
$ws.Range("I5").Value = 1.25
$ws.Range("N5").Value = 1.53
$ws.Range("O5").Value = 2.2
$ws.Range("U5").Value = 50
$ws.Range("V5").Value = 23
$ws.Range("W5").Value = 175
$ws.Range("X5").Value = 80
$ws.Range("Y5").Value = 70
$ws.Range("Z5").Value = 14.5
$ws.Range("AB5").Value = 19
$ws.Range("AC5").Value = 80
$ws.Range("AD5").Value = 500
$ws.Range("AE5").Value = 6.6
$ws.Range("AF5").Value = 5.5
$ws.Range("AG5").Value = 7.7
$ws.Range("AH5").Value = 6.4
$ws.Range("AI5").Value = 9
$ws.Range("AJ5").Value = 22

# Row 8
$ws.Range("J8").Value = 1.04
$ws.Range("L8").Value = 1.25
$ws.Range("P8").Value = 1.33

# Row 9
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 1.7
$ws.Range("L9").Value = 1.22
$ws.Range("M9").Value = 4
$ws.Range("N9").Value = 1.75
$ws.Range("O9").Value = 2.05
$ws.Range("P9").Value = 1.3
$ws.Range("Q9").Value = 3.25
$ws.Range("R9").Value = 1.73
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 15
$ws.Range("U9").Value = 26
$ws.Range("Z9").Value = 12

# Row 10
$ws.Range("G10").Value = 2.3
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 3.1
$ws.Range("J10").Value = 1.05
$ws.Range("K10").Value = 11
$ws.Range("L10").Value = 1.29
$ws.Range("M10").Value = 3.5
$ws.Range("N10").Value = 1.93
$ws.Range("O10").Value = 1.88
$ws.Range("T10").Value = 8
$ws.Range("W10").Value = 21
$ws.Range("X10").Value = 19
$ws.Range("Z10").Value = 10
$ws.Range("AA10").Value = 6.5
$ws.Range("AE10").Value = 9.5
$ws.Range("AH10").Value = 34
$ws.Range("AI10").Value = 26

# Row 12
$ws.Range("J12").Value = 1.03
$ws.Range("L12").Value = 1.22
$ws.Range("P12").Value = 1.3

# Row 13
$ws.Range("L13").Value = 1.29
$ws.Range("M13").Value = 3.5
$ws.Range("N13").Value = 1.9
$ws.Range("O13").Value = 1.9

# Row 14
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 3.5
$ws.Range("I14").Value = 3.5
$ws.Range("Q14").Value = 3.25
$ws.Range("R14").Value = 1.58
$ws.Range("W14").Value = 19
$ws.Range("Y14").Value = 23
$ws.Range("AB14").Value = 12
